# Remove a column from the alcohol measurement data on Sheet1.
# Column M (13) is deleted; the old column N shifts left to become the new M.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Columns("M:M").Delete()

$ws.Range("M1").Select()
